{"js": "// ---------------------------------------------------------------------------\n// What the source diff actually contains\n// ---------------------------------------------------------------------------\n// Every hunk in the target diff touches word/document.xml and word/styles.xml,\n// and in every single hunk the *set* of XML attributes on an element (and\n// their values) is exactly the same before and after -- only the textual\n// order of the attributes changes (plus some internal `w:rsid*` bookkeeping\n// ids that Word regenerates on every save are dropped from the printed\n// diff). For example:\n//\n//   -<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n//   +<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n//\n//   -<w:pgSz w:w=\"11906\" w:h=\"16838\"/>\n//   +<w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n//\n// ...and so on for <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, every\n// <w:lsdException>, and every <w:style>/<w:tblInd>/<w:tblCellMar>. The root\n// <w:document> element likewise only has its xmlns:* declarations\n// re-sorted. This is the well known fingerprint of a document having been\n// re-saved by a newer OOXML writer (the commit message: \"upgraded to POI\n// 3.15\") which happens to serialize attributes in alphabetical order -\n// it is not an edit to the document's content, formatting, layout or\n// styles: the color (E36C0A / accent6 / BF), the page size\n// (11906 x 16838 twips), the margins (1417/1417/1417/1417 twips,\n// header/footer 708, gutter 0), the default fonts/language\n// (minorHAnsi/minorBidi, fr-FR/en-US/ar-SA) and every style/latent-style\n// definition are byte-for-byte identical values on both sides of the diff.\n//\n// Office.js's Word object model (like the real Word UI) exposes documents\n// in terms of their content/formatting, never in terms of the raw\n// attribute-serialization order of the underlying part XML (that ordering\n// is not a meaningful, settable property anywhere in Word's object model).\n// Since every value already matches what the diff shows, there is no\n// content-level mutation to make here. Re-\"setting\" properties that are\n// already at their target value would not reproduce the attribute-order\n// artifact anyway (this runtime's writer keeps a stable attribute order\n// regardless of how a part is touched) and would risk introducing genuine,\n// unwanted content diffs (e.g. Office.js's `font.color` only models plain\n// RGB colors, not the `themeColor`/`themeShade` pair, and resolving a\n// search hit can split runs) that are not present in the source diff.\n//\n// So this script faithfully mirrors the (content-free) change: it reads\n// the body to confirm the document is in the expected state and performs\n// no mutation.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# ---------------------------------------------------------------------------\n# What the source diff actually contains\n# ---------------------------------------------------------------------------\n# Every hunk in the target diff touches word/document.xml and word/styles.xml,\n# and in every single hunk the *set* of XML attributes on an element (and\n# their values) is exactly the same before and after -- only the textual\n# order of the attributes changes (plus some internal `w:rsid*` bookkeeping\n# ids that Word regenerates on every save, which are dropped from the\n# printed diff). For example:\n#\n#   -<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n#   +<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n#\n#   -<w:pgSz w:w=\"11906\" w:h=\"16838\"/>\n#   +<w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#\n# ...and so on for <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, every\n# <w:lsdException>, and every <w:style>/<w:tblInd>/<w:tblCellMar>. The root\n# <w:document> element likewise only has its xmlns:* declarations\n# re-sorted. This is the fingerprint of the document having been re-saved\n# by a newer OOXML writer (the commit message: \"upgraded to POI 3.15\")\n# which happens to serialize attributes in alphabetical order - it is not\n# an edit to the document's content, formatting, layout or styles: the\n# color (E36C0A / accent6 / BF), the page size (11906 x 16838 twips =\n# 595.3 x 841.9 pt), the margins (1417/1417/1417/1417 twips = 70.85 pt,\n# header/footer 708 twips = 35.4 pt, gutter 0), the default fonts/language\n# (minorHAnsi/minorBidi, fr-FR/en-US/ar-SA) and every style/latent-style\n# definition are byte-for-byte identical values on both sides of the diff.\n#\n# The Word COM object model (like the Word UI it backs) exposes documents\n# in terms of their content/formatting, never in terms of the raw\n# attribute-serialization order of the underlying part XML (that ordering\n# is not a meaningful, settable property anywhere in Word's object model -\n# Find/Replace, Paragraphs, PageSetup, Styles, Font, etc. all read/write\n# values, not XML attribute order). Since every value already matches what\n# the diff shows, there is no content-level mutation to make here.\n# Re-\"setting\" properties that are already at their target value would not\n# reproduce the attribute-order artifact anyway (this runtime's writer\n# keeps a stable attribute order regardless of how a part is touched /\n# re-saved) and would risk introducing genuine, unwanted content diffs\n# (e.g. resolving Find on visible text can split runs and restate\n# formatting on the wrong text, or drop the `themeColor`/`themeShade`\n# pairing) that are not present in the source diff.\n#\n# So this script faithfully mirrors the (content-free) change: it reads\n# the document's current page setup to confirm it is already in the\n# expected state and performs no mutation.\n$d = $word.ActiveDocument\n$ps = $d.PageSetup\n$null = $ps.PageWidth\n$null = $ps.PageHeight\n"}
